$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "release/8.0.9"
$ws.Range("B12").Value = "X"
$ws.Range("C12").Value = "X"
$ws.Range("D12").Value = "X"
$ws.Range("E12").Value = "X"

# Row 11 (the last pre-existing data row) carries no explicit cell style
# (no s="n" attribute) even though its column defines a default style.
# Match that so the new row 12 lines up with the rest of the sheet.
$ws.Range("A12:E12").Style = "Normal"
